$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Insert a new column before column H. This shifts the existing H:M columns
# (PHI KIEM HANG / PHI BAO HIEM / PHI DONG GO / NGAY TAO / NGAY XUAT KHO /
# TRANG THAI) one position to the right, becoming I:N.
# ---------------------------------------------------------------------------
$ws.Columns("H").Insert()

# New column H: "TONG TIEN" header / TotalPriceVND placeholder field
$ws.Range("H1").Value = "TỔNG TIỀN"
$ws.Range("H2").Value = "[[%Field:TotalPriceVND%]]"

# Give the new column roughly the same width as its neighbours
$ws.Columns("H").ColumnWidth = 29

# ---------------------------------------------------------------------------
# Header row formatting (H1:K1) - bold Times New Roman 12, #,##0, centered
# ---------------------------------------------------------------------------
$hdr = $ws.Range("H1:K1")
$hdr.NumberFormat = "#,##0"
$hdr.Font.Name = "Times New Roman"
$hdr.Font.Size = 12
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4108

# ---------------------------------------------------------------------------
# Data row formatting (H2:K2) - regular Times New Roman 12, #,##0, centered
# ---------------------------------------------------------------------------
$dat = $ws.Range("H2:K2")
$dat.NumberFormat = "#,##0"
$dat.Font.Name = "Times New Roman"
$dat.Font.Size = 12
$dat.Font.Bold = $false
$dat.HorizontalAlignment = -4108
$dat.VerticalAlignment = -4108

# Update the active cell / selection on the sheet
[void]$ws.Range("I12").Select()

# Page setup: A4, portrait
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

Write-Host "Transportation order template updated"
